# Changing to work against generic APIs
#
# 1) Rename the "config" sheet to "All".
# 2) Insert a brand-new "TestAll" sheet in front of everything else, with a
#    couple of sample rows that hit a generic public API (jsonplaceholder)
#    instead of the old internal-only domain.
# 3) Point the selections at the cells the author left active, and re-freeze
#    the header row on the new sheet just like the existing ones.

$wb = $excel.ActiveWorkbook

# --- Rename "config" -> "All" ---------------------------------------------
$wb.Worksheets.Item("config").Name = "All"

# --- New "TestAll" sheet, inserted as the first tab ------------------------
# NOTE: inserting a sheet shifts the tab position of every sheet after it,
# so any worksheet reference grabbed *before* this Add() becomes stale.
# Re-fetch LoadParms/All by name below, only after the insert has happened.
$testAll = $wb.Worksheets.Add($wb.Worksheets.Item(1))
$testAll.Name = "TestAll"

# Header row - reuses the same column headings as LoadParms/All.
$testAll.Cells.Item(1,1).Value = "ID"
$testAll.Cells.Item(1,2).Value = "Active"
$testAll.Cells.Item(1,3).Value = "DependentID"
$testAll.Cells.Item(1,4).Value = "FieldToLoad"
$testAll.Cells.Item(1,5).Value = "SourceField"
$testAll.Cells.Item(1,6).Value = "Domain"
$testAll.Cells.Item(1,7).Value = "api"
$testAll.Cells.Item(1,8).Value = "Parms"

# Row 2 - GET /posts (list), keyed by userId.
$testAll.Cells.Item(2,1).Value = 1
$testAll.Cells.Item(2,2).Value = "Y"
$testAll.Cells.Item(2,6).Value = "https://jsonplaceholder.typicode.com"
$testAll.Cells.Item(2,7).Value = "/posts"
$testAll.Cells.Item(2,4).Value = "userId"
$testAll.Cells.Item(2,5).Value = "userId"

# Row 3 - GET /posts/{userId} (single post), depends on row 2's ID.
$testAll.Cells.Item(3,1).Value = 2
$testAll.Cells.Item(3,2).Value = "Y"
$testAll.Cells.Item(3,3).Value = 1
$testAll.Cells.Item(3,6).Value = "https://jsonplaceholder.typicode.com"
$testAll.Cells.Item(3,7).Value = "/posts/{userId}"

# Column widths, approximated to match the author's auto-fit sizing.
$testAll.Columns.Item(1).ColumnWidth = 1.83
$testAll.Columns.Item(2).ColumnWidth = 5.17
$testAll.Columns.Item(3).ColumnWidth = 11
$testAll.Columns.Item(4).ColumnWidth = 10.17
$testAll.Columns.Item(5).ColumnWidth = 9.67
$testAll.Columns.Item(6).ColumnWidth = 32.17
$testAll.Columns.Item(7).ColumnWidth = 12.83
$testAll.Columns.Item(8).ColumnWidth = 16.33
$testAll.Columns.Item(9).ColumnWidth = 2.33
$testAll.Columns.Item(10).ColumnWidth = 42.97
$testAll.Columns.Item(13).ColumnWidth = 137.17

# Freeze the header row, then leave the selection where the author left it.
$testAll.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$testAll.Range("F10").Select()

# --- Restore the selections left on the other two sheets -------------------
# Fetched fresh (by name) now that the tab order has settled.
$loadParms = $wb.Worksheets.Item("LoadParms")
$loadParms.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$loadParms.Range("D23").Select()

$allSheet = $wb.Worksheets.Item("All")
$allSheet.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$allSheet.Range("D8").Select()

# Leave "TestAll" as the active, selected tab.
$testAll.Select()
$testAll.Range("F10").Select()
